# Abstract updates. Overview figure update.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The main diagram lives inside a top-level group shape ("Group 156").
$grp = $s.Shapes.Item("Group 156")

# --- Change 1: "TreeBeard" -> "Treebeard" --------------------------------
$treeBeardShape = $grp.GroupItems.Item("Rectangle 14")
$treeBeardShape.TextFrame.TextRange.Text = "Treebeard"

# --- Change 2: "Intel`t`tAMD`tArm" -> "x86`t`t" + "Risc" + "-V" + "`tArm" --
$hwShape = $grp.GroupItems.Item("Rectangle 48")
$hwRange = $hwShape.TextFrame.TextRange
$hwRange.Text = "x86`t`tRisc-V`tArm"
# Split the single run into 4 runs matching the new wording, without
# touching any formatting (re-assigning the same text on a sub-range is
# enough to force a run break at that boundary).
$r2 = $hwRange.Characters(6, 4)
$r2.Text = $r2.Text
$r3 = $hwRange.Characters(10, 2)
$r3.Text = $r3.Text
$r4 = $hwRange.Characters(12, 4)
$r4.Text = $r4.Text

# --- Change 3: "vectorization" -> "Vectorization" ------------------------
$vecShape = $grp.GroupItems.Item("Rectangle 145")
$vecShape.TextFrame.TextRange.Text = "Vectorization"

# --- Change 4: fix "LIghtGBM" -> "LightGBM" -------------------------------
$gbmShape = $grp.GroupItems.Item("Rectangle 155")
$gbmRange = $gbmShape.TextFrame.TextRange
# "XGBoost" (7) + "      " (6) + "LIghtGBM" (8) -> the misspelled run starts
# at character 14.
$gbmSub = $gbmRange.Characters(14, 8)
$gbmSub.Text = "LightGBM"

# --- Change 5: "parallelize" -> "Parallelize" -----------------------------
$parallelShape = $s.Shapes.Item("Rectangle 2")
$parallelShape.TextFrame.TextRange.Text = "Parallelize"
